$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.34021225948807
$ws.Range("C2").Value = 11.32962030432695
$ws.Range("E2").Value = 13.4550518476423
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 30.09752767458116
$ws.Range("H2").Value = 14.793806425028
$ws.Range("K2").Value = 7.750541073583955
$ws.Range("L2").Value = 9.888356956101726
$ws.Range("M2").Value = 13.61510343499229
$ws.Range("N2").Value = 19.54999094260289
$ws.Range("O2").Value = 22.66362574558043
$ws.Range("B3").Value = 11.12324507180015
$ws.Range("C3").Value = 11.35220880986481
$ws.Range("E3").Value = 13.48645443091365
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 30.22426053232374
$ws.Range("H3").Value = 14.84081388119339
$ws.Range("K3").Value = 7.582559042601105
$ws.Range("L3").Value = 9.89432869703689
$ws.Range("M3").Value = 13.58291167343278
$ws.Range("N3").Value = 19.60262177698279
$ws.Range("O3").Value = 22.74867085438285
$ws.Range("B4").Value = 10.98961040411148
$ws.Range("C4").Value = 11.36694878635508
$ws.Range("E4").Value = 13.50768776459958
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 30.30970854198389
$ws.Range("H4").Value = 14.87156190725347
$ws.Range("K4").Value = 7.478340101984481
$ws.Range("L4").Value = 9.899278911124268
$ws.Range("M4").Value = 13.56489077449253
$ws.Range("N4").Value = 19.63657658063364
$ws.Range("O4").Value = 22.80474599032967
$ws.Range("B5").Value = 10.93512183525843
$ws.Range("C5").Value = 11.37317491185382
$ws.Range("E5").Value = 13.51683162801384
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 30.34644378372098
$ws.Range("H5").Value = 14.88456663355644
$ws.Range("K5").Value = 7.435655946782386
$ws.Range("L5").Value = 9.901619622035051
$ws.Range("M5").Value = 13.55799114790988
$ws.Range("N5").Value = 19.65082671151991
$ws.Range("O5").Value = 22.82856694231684
$ws.Range("B6").Value = 10.92607432674154
$ws.Range("C6").Value = 11.37422202693017
$ws.Range("E6").Value = 13.51837962831895
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 30.35265912375429
$ws.Range("H6").Value = 14.88675474409678
$ws.Range("K6").Value = 7.428557045338708
$ws.Range("L6").Value = 9.902027850603725
$ws.Range("M6").Value = 13.55687243599557
$ws.Range("N6").Value = 19.65321792704756
$ws.Range("O6").Value = 22.83258097132548
$ws.Range("B7").Value = 10.98887558126154
$ws.Range("C7").Value = 11.36703186470336
$ws.Range("E7").Value = 13.50780909313322
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 30.31019622134151
$ws.Range("H7").Value = 14.87173537088121
$ws.Range("K7").Value = 7.477765238805325
$ws.Range("L7").Value = 9.899309168217259
$ws.Range("M7").Value = 13.56479591913087
$ws.Range("N7").Value = 19.63676708813384
$ws.Range("O7").Value = 22.80506332102497
$ws.Range("B8").Value = 11.26553171297321
$ws.Range("C8").Value = 11.33722847695532
$ws.Range("E8").Value = 13.46547447725666
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 30.13963783836934
$ws.Range("H8").Value = 14.80962366690818
$ws.Range("K8").Value = 7.692877872300873
$ws.Range("L8").Value = 9.890150052411482
$ws.Range("M8").Value = 13.60364478591942
$ws.Range("N8").Value = 19.56779847058596
$ws.Range("O8").Value = 22.69214874764325
$ws.Range("B9").Value = 11.80163083791589
$ws.Range("C9").Value = 11.28566623721214
$ws.Range("E9").Value = 13.39793299747578
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 29.86596217170812
$ws.Range("H9").Value = 14.70275464345141
$ws.Range("K9").Value = 8.103785188070265
$ws.Range("L9").Value = 9.882340727553672
$ws.Range("M9").Value = 13.69343013234871
$ws.Range("N9").Value = 19.44550879453419
$ws.Range("O9").Value = 22.50132861034478
$ws.Range("B10").Value = 12.18757443620318
$ws.Range("C10").Value = 11.2519440822221
$ws.Range("E10").Value = 13.3577283591835
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 29.70224528409369
$ws.Range("H10").Value = 14.63330356304789
$ws.Range("K10").Value = 8.39601285125652
$ws.Range("L10").Value = 9.882745365477742
$ws.Range("M10").Value = 13.76734548632091
$ws.Range("N10").Value = 19.36349229397878
$ws.Range("O10").Value = 22.37979132968457
$ws.Range("B11").Value = 12.36063670390971
$ws.Range("C11").Value = 11.23749906490138
$ws.Range("E11").Value = 13.34147920372711
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 29.63593735714563
$ws.Range("H11").Value = 14.60366925046587
$ws.Range("K11").Value = 8.526277596857598
$ws.Range("L11").Value = 9.884251819063358
$ws.Range("M11").Value = 13.80261876070499
$ws.Range("N11").Value = 19.32786626909463
$ws.Range("O11").Value = 22.32855287477526
$ws.Range("B12").Value = 12.42574908636362
$ws.Range("C12").Value = 11.23215730140395
$ws.Range("E12").Value = 13.33561906524219
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 29.61200755867572
$ws.Range("H12").Value = 14.59272871250515
$ws.Range("K12").Value = 8.575177295740396
$ws.Range("L12").Value = 9.885011328610078
$ws.Range("M12").Value = 13.81620576921196
$ws.Range("N12").Value = 19.31461665107308
$ws.Range("O12").Value = 22.30973268358599
$ws.Range("B13").Value = 12.4117457153051
$ws.Range("C13").Value = 11.23330204967287
$ws.Range("E13").Value = 13.33686812202464
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 29.61710872829209
$ws.Range("H13").Value = 14.59507244724384
$ws.Range("K13").Value = 8.564665615679619
$ws.Range("L13").Value = 9.88483936514792
$ws.Range("M13").Value = 13.81326946325653
$ws.Range("N13").Value = 19.3174594812887
$ws.Range("O13").Value = 22.31376002639161
$ws.Range("B14").Value = 12.36600235094369
$ws.Range("C14").Value = 11.23705702731166
$ws.Range("E14").Value = 13.34099121503085
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 29.63394496941179
$ws.Range("H14").Value = 14.60276352983095
$ws.Range("K14").Value = 8.530309441364242
$ws.Range("L14").Value = 9.884310522187942
$ws.Range("M14").Value = 13.80373200964149
$ws.Range("N14").Value = 19.32677138662388
$ws.Range("O14").Value = 22.32699284585764
$ws.Range("B15").Value = 12.33792640579881
$ws.Range("C15").Value = 11.23937374683708
$ws.Range("E15").Value = 13.34355488452554
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 29.64441141596487
$ws.Range("H15").Value = 14.60751116571227
$ws.Range("K15").Value = 8.509208188255716
$ws.Range("L15").Value = 9.884011176042172
$ws.Range("M15").Value = 13.79791974725934
$ws.Range("N15").Value = 19.33250657950919
$ws.Range("O15").Value = 22.33517422642738
$ws.Range("B16").Value = 12.17620858965222
$ws.Range("C16").Value = 11.2529060721064
$ws.Range("E16").Value = 13.3588313071294
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 29.7067434555293
$ws.Range("H16").Value = 14.63527962595162
$ws.Range("K16").Value = 8.387442106968427
$ws.Range("L16").Value = 9.882673431901821
$ws.Range("M16").Value = 13.76507285412756
$ws.Range("N16").Value = 19.36585433902293
$ws.Range("O16").Value = 22.38322140194683
$ws.Range("B17").Value = 12.07631280943392
$ws.Range("C17").Value = 11.2614366808259
$ws.Range("E17").Value = 13.36872521864165
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 29.74707775082602
$ws.Range("H17").Value = 14.65281620854823
$ws.Range("K17").Value = 8.312025669641598
$ws.Range("L17").Value = 9.882190698771378
$ws.Range("M17").Value = 13.74533953067183
$ws.Range("N17").Value = 19.38674267867842
$ws.Range("O17").Value = 22.41373422696612
$ws.Range("B18").Value = 12.0186230778384
$ws.Range("C18").Value = 11.26642756650998
$ws.Range("E18").Value = 13.374607982423
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 29.77104546168364
$ws.Range("H18").Value = 14.6630872327651
$ws.Range("K18").Value = 8.26839918318197
$ws.Range("L18").Value = 9.882037583519915
$ws.Range("M18").Value = 13.73414509341077
$ws.Range("N18").Value = 19.39891563098315
$ws.Range("O18").Value = 22.43166552724204
$ws.Range("B19").Value = 11.99905238082542
$ws.Range("C19").Value = 11.26813188947122
$ws.Range("E19").Value = 13.37663277827253
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 29.77929234849215
$ws.Range("H19").Value = 14.66659651671626
$ws.Range("K19").Value = 8.253586659561869
$ws.Range("L19").Value = 9.882007163709751
$ws.Range("M19").Value = 13.7303818005272
$ws.Range("N19").Value = 19.40306443997047
$ws.Range("O19").Value = 22.43780220521704
$ws.Range("B20").Value = 12.08697139262368
$ws.Range("C20").Value = 11.26051986134114
$ws.Range("E20").Value = 13.36765212078105
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 29.74270452302016
$ws.Range("H20").Value = 14.65093032300841
$ws.Range("K20").Value = 8.320079976523076
$ws.Range("L20").Value = 9.8822292048936
$ws.Range("M20").Value = 13.74742411974083
$ws.Range("N20").Value = 19.38450267858058
$ws.Range("O20").Value = 22.41044663629727
$ws.Range("B21").Value = 12.3794502318213
$ws.Range("C21").Value = 11.23595062262668
$ws.Range("E21").Value = 13.33977221159002
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 29.62896770447962
$ws.Range("H21").Value = 14.60049684039976
$ws.Range("K21").Value = 8.540412658265788
$ws.Range("L21").Value = 9.884460734689249
$ws.Range("M21").Value = 13.80652721421814
$ws.Range("N21").Value = 19.32402971596379
$ws.Range("O21").Value = 22.32309022852172
$ws.Range("B22").Value = 12.56811285909376
$ws.Range("C22").Value = 11.22064052122463
$ws.Range("E22").Value = 13.3232591158006
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 29.56151180110435
$ws.Range("H22").Value = 14.56917528401641
$ws.Range("K22").Value = 8.681895347611556
$ws.Range("L22").Value = 9.88702051723541
$ws.Range("M22").Value = 13.84649029484054
$ws.Range("N22").Value = 19.2859125246342
$ws.Range("O22").Value = 22.26939437696759
$ws.Range("B23").Value = 12.46766745534674
$ws.Range("C23").Value = 11.22874359515057
$ws.Range("E23").Value = 13.33191628963915
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 29.59688335103118
$ws.Range("H23").Value = 14.58574229535131
$ws.Range("K23").Value = 8.606627648470006
$ws.Range("L23").Value = 9.885553925951667
$ws.Range("M23").Value = 13.82504150336654
$ws.Range("N23").Value = 19.3061280976057
$ws.Range("O23").Value = 22.29774194045895
$ws.Range("B24").Value = 12.08215344698595
$ws.Range("C24").Value = 11.26093408588207
$ws.Range("E24").Value = 13.36813666208222
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 29.74467923322542
$ws.Range("H24").Value = 14.65178234307988
$ws.Range("K24").Value = 8.316439457495147
$ws.Range("L24").Value = 9.882211408659693
$ws.Range("M24").Value = 13.74648120764046
$ws.Range("N24").Value = 19.38551487180618
$ws.Range("O24").Value = 22.41193174414526
$ws.Range("B25").Value = 11.65771456075827
$ws.Range("C25").Value = 11.29888203708117
$ws.Range("E25").Value = 13.41454951449291
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 29.93345992741104
$ws.Range("H25").Value = 14.73007088876406
$ws.Range("K25").Value = 7.994114890330914
$ws.Range("L25").Value = 9.883371067021125
$ws.Range("M25").Value = 13.667718695334
$ws.Range("N25").Value = 19.47721143240253
$ws.Range("O25").Value = 22.54967403599906
